$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.031337738037109
$ws.Range("B1").Value = 3.286305665969849
$ws.Range("C1").Value = 3.620364189147949
$ws.Range("D1").Value = 2.017948627471924
$ws.Range("E1").Value = 1.173141121864319
